# Workbook-level edit script.
#
# Commit message: "#2 #3 display and delete comment"
#   #2 "display" -> the POSTS table on the "Sheet2" worksheet was missing the
#                    FK/userid/int row (it showed a blank row where the FK
#                    relationship should be); it is now shown.
#   #3 "delete"   -> the "likes" table on the "Sheet3" worksheet had a stray
#                    duplicate FK/accountid/int row; it is removed.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet2: "POSTS" table (H14:J19) -------------------------------------
# Before:
#   H15:J15 = PK / postId   / int
#   H16:J16 = (blank row)
#   H17:J17 = FK / userid   / int
#   H18:J18 = (blank) / description / varchar(50)
#   H19:J19 = (blank) / datePost    / date
#   H20:J20 = (blank row)
#
# After: the blank separator row is removed and every row below it shifts
# up by one, so the FK row becomes visible right under the PK row, and the
# table ends one row earlier (dimension shrinks from P20 to P19).

$ws2.Rows(20).Delete()

$ws2.Range("H16").Value = "FK"
$ws2.Range("I16").Value = "userid"
$ws2.Range("J16").Value = "int"

$ws2.Range("H17").ClearContents()
$ws2.Range("I17").Value = "description"
$ws2.Range("J17").Value = "varchar(50)"

$ws2.Range("H18").ClearContents()
$ws2.Range("I18").Value = "datePost"
$ws2.Range("J18").Value = "date"

$ws2.Range("H19:J19").Clear()

# --- Sheet3: "likes" table (M3:O6) ----------------------------------------
# Before:
#   M5:O5 = FK / accountid / int   <- stray duplicate row, deleted
#   M6:O6 = FK / userId    / int
#   M7:O7 = FK / postid    / int
#
# After: row 5 is removed and rows 6-7 shift up by one.

$ws3.Range("N5").Value = "userId"
$ws3.Range("M6").Value = "FK"
$ws3.Range("N6").Value = "postid"
$ws3.Range("O6").Value = "int"
$ws3.Range("M7:O7").Clear()

# --- View state: Sheet3 becomes the active tab, matching the saved file ---
$ws2.Activate()
$ws2.Range("J23").Select()

$ws3.Activate()
$ws3.Range("M4").Select()
